# DPLKAKT068-001 - Setup Periode Bulanan - "Update Regresi Tanggal 31/03/2023"
# Roll the regression period forward one year: 2023/04 -> 2024/04.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TGL_AWAL (O2): 15/04/2023 -> 15/04/2024
$ws.Range("O2").Value = "15/04/2024"

# PERIODE_BULANAN (Q2): 202305 -> 202405
$ws.Range("Q2").Value = "202405"

# VERIFIKASI (T2): stays "05" (kept explicit since the underlying shared-string slot moves)
$ws.Range("T2").Value = "05"

# Reflect the saved view state: scrolled right so column O is the leftmost visible
# column, with the active selection on X2.
$excel.ActiveWindow.ScrollColumn = 15
[void]$ws.Range("X2").Select()
